# Update Excel odds values for rows 3 and 4 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 changes
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 2.1
$ws.Range("K3").Value = 2.1
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 6
$ws.Range("Z3").Value = 9.5
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 23
$ws.Range("AL3").Value = 51
$ws.Range("AW3").Value = 8
$ws.Range("AY3").Value = 41

# Row 4 changes
$ws.Range("I4").Value = 2.42
$ws.Range("J4").Value = 3.35
$ws.Range("L4").Value = 2.92
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 1.91
$ws.Range("R4").Value = 1.8
$ws.Range("W4").Value = 8.75
$ws.Range("X4").Value = 14
$ws.Range("Y4").Value = 10.25
$ws.Range("AA4").Value = 24
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 9.5
$ws.Range("AD4").Value = 6.1
$ws.Range("AE4").Value = 13
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 13
$ws.Range("AL4").Value = 19.5
$ws.Range("AM4").Value = 27
$ws.Range("AN4").Value = 4.75
$ws.Range("AO4").Value = 15
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 70
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.62
$ws.Range("AU4").Value = 6.5
$ws.Range("AW4").Value = 4.4
$ws.Range("AY4").Value = 18
$ws.Range("BA4").Value = 70
